$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that sat between "Soil; " and
#    "Accessibility, only to transform into normal soil, " in the
#    "Create" Statistics line.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Add new bookmarks (OLE_LINK67, OLE_LINK68, and a new _GoBack) that
#    wrap the whole paragraph containing "You throw a ball of stone...".
$rng = $d.Content.Duplicate
$rng.Find.Execute("You throw a ball of stone", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range

# Add in reverse order so the resulting bookmarkStart order (by id) comes
# out as OLE_LINK67, OLE_LINK68, _GoBack -- matching the target document.
$d.Bookmarks.Add("_GoBack", $para)
$d.Bookmarks.Add("OLE_LINK68", $para)
$d.Bookmarks.Add("OLE_LINK67", $para)
